$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$tbl = $ws.ListObjects.Item("Table1")

# --- 1. Insert a new row at row 98 (shifts old rows 98-140 down to 99-141) ---
$ws.Rows("98:98").Insert()

# --- 2. Resize the table to include the new row (A8:K140 -> A8:K141) ---
$tbl.Resize($ws.Range("A8:K141"))

# --- 3. Copy formatting from row 81 (the "2023" year-marker row) into the new row 98 ---
$ws.Range("A81:K81").Copy()
$ws.Range("A98:K98").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 4. Populate the new row 98 as the "2024" year marker ---
$ws.Range("A98").Value = "'2024"
$ws.Range("G98").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"

# --- 5. Row 94: add EARNED value ---
$ws.Range("C94").Value = 1.25

# --- 6. Row 95: SL leave taken, earned, absence w/ pay hour, and return date ---
$ws.Range("B95").Value = "SL(1-0-0)"
$ws.Range("C95").Value = 1.25
$ws.Range("H95").Value = 1
$ws.Range("K93").Copy()
$ws.Range("K95").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("K95").Value = 45212

# --- 7. Row 96: SL leave taken, earned, absence w/ pay hour, and return date ---
$ws.Range("B96").Value = "SL(1-0-0)"
$ws.Range("C96").Value = 1.25
$ws.Range("H96").Value = 1
$ws.Range("K93").Copy()
$ws.Range("K96").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("K96").Value = 45260

# --- 8. Row 97: VL leave taken, absence w/o pay days, and remarks ---
$ws.Range("B97").Value = "VL(2-0-0)"
$ws.Range("D97").Value = 2
$ws.Range("K97").Value = "12/19,20/2023"

# --- 9. The table resize re-wrote the last row's calculated-column formula using
#        the "[@EARNED]" shorthand, which this engine evaluates to #VALUE! for a
#        blank input; restore the normal structured-reference formula so it
#        evaluates back to the expected blank string, matching every other row.
$ws.Range("G141").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"

# --- 10. Put the selection where the author last left it ---
$ws.Range("K97").Select()
